$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.075.91"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.929.15"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'326.36"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "'0.4608"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").Value = "'0.3830"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "'0.07763"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "'0.9795"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "'22.55"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "1.972.63"
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("D13").Value = "'6.982"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'5.698"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'0.07063"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D17").Value = "'84.43"
$ws.Range("E17").Value = "  +2.43%  "
$ws.Range("D18").Value = "'0.000009541"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "'16.76"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "29.116.72"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "'5.347"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'10.96"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "'2.078"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").Value = "'158.01"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "'19.13"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").Value = "'5.670"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "'118.15"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "'1.850"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("D30").Value = "'0.09343"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "'0.8652"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("D32").Value = "'5.139"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").Value = "'1.250"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").Value = "'3.011"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'0.05707"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "'1.004"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'0.02053"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").Value = "'3.070"
$ws.Range("E39").Value = "  +13.85%  "
$ws.Range("D40").Value = "'7.548"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").Value = "'0.1757"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'9.373"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "'0.000002822"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").Value = "'2.203"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("D46").Value = "'0.5207"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'11.27"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.06928"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "'1.780"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "'110.36"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.31%  "
